$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Update simple property values that changed in place
$ws.Range("B3").Value = "6.0.0"
$ws.Range("B8").Value = "2022-01-21T20:46:54+00:00"
$ws.Range("B9").Value = "Alvearie Team"

# Remove the duplicated second "Contact" row (old row 11); this shifts
# everything below it up by one row.
$ws.Rows.Item(11).Delete()

# The old "Contact" row (now row 10) becomes the new "Jurisdiction" row.
$ws.Range("A10").Value = "Jurisdiction"
$ws.Range("B10").Value = "United States of America"

# "Case Sensitive" is now row 14; give it a value of "true".
$ws.Range("B14").Value = "true"
